# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts.
#
# This updates the DAMSLTag (column I) and DialogAct (column J) values
# for a set of rows on Sheet1, reflecting the re-annotated dialog acts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=10;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=17;  I="sv"; J="Statement-opinion"},
    @{Row=49;  I="sd"; J="Statement-non-opinion"},
    @{Row=53;  I="sd"; J="Statement-non-opinion"},
    @{Row=58;  I="sd"; J="Statement-non-opinion"},
    @{Row=65;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=69;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=70;  I="sv"; J="Statement-opinion"},
    @{Row=78;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=79;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=104; I="sd"; J="Statement-non-opinion"},
    @{Row=105; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=111; I="sv"; J="Statement-opinion"},
    @{Row=114; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=119; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=179; I="aa"; J="Agree/Accept"},
    @{Row=185; I="aa"; J="Agree/Accept"},
    @{Row=189; I="aa"; J="Agree/Accept"},
    @{Row=190; I="sv"; J="Statement-opinion"},
    @{Row=196; I="sd"; J="Statement-non-opinion"},
    @{Row=208; I="ba"; J="Appreciation"},
    @{Row=216; I="sv"; J="Statement-opinion"},
    @{Row=225; I="aa"; J="Agree/Accept"},
    @{Row=232; I="sd"; J="Statement-non-opinion"},
    @{Row=235; I="%";  J="Uninterpretable"},
    @{Row=239; I="aa"; J="Agree/Accept"},
    @{Row=240; I="sv"; J="Statement-opinion"},
    @{Row=245; I="sv"; J="Statement-opinion"},
    @{Row=257; I="sd"; J="Statement-non-opinion"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value  = $u.I   # Column I = DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.J   # Column J = DialogAct
}
